# neptunes-champions (Version 2): add a "Meta description" paragraph under
# the title, drop the duplicated bold title paragraph at the bottom of the
# doc, and turn the leftover italic "meta description" paragraph there into
# an AI image-generation prompt for the article's feature image.

$d = $word.ActiveDocument

# --- Step 1: remove the duplicate bold title paragraph near the bottom ---
# (it is identical to paragraph 1's text, so skip paragraph 1 itself)
$n = $d.Paragraphs.Count
for ($i = $n; $i -ge 2; $i--) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text.StartsWith("Play Neptune's Champions for Free: Review and Ratings")) {
    $p.Range.Delete()
    break
  }
}

# --- Step 2: turn the trailing italic paragraph into the image prompt ---
# Do this before inserting the new "Meta description" paragraph up top,
# since that insertion will introduce a second (non-italic) copy of this
# same sentence and we only want to touch the original, italic, one.
$oldText = "Read our review of Neptune's Champions, a medium volatility slot with a fun free spins feature and potential for big wins. Play it for free!"
$newText = "Create a cartoon-style feature image for Neptune's Champions online slot game. The image should feature a happy Maya warrior with glasses. The warrior should be standing in front of a background of ocean waves, with Neptune looming in the distance. Make sure to include the title of the game in the image, along with any other relevant symbols or logos. The overall feel of the image should be fun and exciting, inviting players to dive into the game and discover its treasures."

$rng = $d.Content
$found = $rng.Find.Execute($oldText)
if ($found) {
  # Assign .Text directly (rather than passing $newText as Find's Replace
  # argument) so straight apostrophes aren't auto-corrected into curly
  # ones; this also preserves the run's existing <w:i/> formatting.
  $rng.Text = $newText
}

# --- Step 3: insert the "Meta description" paragraph right after the title ---
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$pStart = $metaPara.Range.Start
$pEnd = $metaPara.Range.End
$target = $d.Range($pStart, $pEnd)

# Build the new paragraph's exact run structure (leading empty run, bold
# "Meta description" run, plain run with the rest of the sentence) via
# InsertXML so the three runs come out exactly as wanted in one shot.
$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Neptune''s Champions, a medium volatility slot with a fun free spins feature and potential for big wins. Play it for free!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($metaXml) | Out-Null
